$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Status" header in column J, matching the existing header style
# by copying the formatting from the adjacent header cell (I1) before writing
# the value.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Status"

# Resize columns B-I slightly (as a result of the new column being added /
# content re-fitted) and set the width of the new column J.
$ws.Columns.Item(2).ColumnWidth = 18.6666666667
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 19
$ws.Columns.Item(7).ColumnWidth = 17
$ws.Columns.Item(8).ColumnWidth = 10.8333333333
$ws.Columns.Item(9).ColumnWidth = 17.3333333333
$ws.Columns.Item(10).ColumnWidth = 26.3333333333

# Move the active selection to D3
$ws.Range("D3").Select() | Out-Null
